$d = $word.ActiveDocument

# --- 1. Title: "Møtereferat-19" -> "Møtereferat-20" ---
$d.Content.Find.Execute("19", $true, $false, $false, $false, $false, $true, 1, $false, "20", 2) | Out-Null

# --- 2. Date: "Dato: 21.04.2021" -> "Dato: 22.04.2021" ---
$d.Content.Find.Execute("Dato: 21", $true, $false, $false, $false, $false, $true, 1, $false, "Dato: 22", 2) | Out-Null

# --- 3. "Gjennomgang av ukens Pull Requests" -> "Gruppen klarer å klarer kollisjonsproblemet." ---
$d.Content.Find.Execute("Gjennomgang av ukens Pull Requests", $true, $false, $false, $false, $false, $true, 1, $false, "Gruppen klarer å klarer kollisjonsproblemet.", 2) | Out-Null

# --- 4. "Testing av spill" -> "Starter med GUI for valg av kort." ---
$d.Content.Find.Execute("Testing av spill", $true, $false, $false, $false, $false, $true, 1, $false, "Starter med GUI for valg av kort.", 2) | Out-Null

# --- 5. Rewrite the "Mathias (Kollisjoner)..." paragraph entirely ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Mathias (Kollisjoner)")) {
        $r = $p.Range
        $full = $d.Range($r.Start, $r.End - 1)
        $full.Text = "Planlegger ny lagmiddag. Og teleskopsaubesøk på Syltøy."
        break
    }
}

# --- 6. Delete the now-obsolete bullet paragraphs ---
# "Gruppen er fornøyd med arbeidet som er gjort til nå"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Gruppen er forn")) {
        $p.Range.Delete()
        break
    }
}

# "Lager Flere arbeidspunkter som må jobbes med"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Lager Flere arbeidspunkter")) {
        $p.Range.Delete()
        break
    }
}

# empty paragraph with ind left=360 directly after the rewritten "Planlegger..." paragraph
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Planlegger ny lagmiddag")) {
        $next = $d.Paragraphs.Item($i + 1)
        if ($next.Range.Text.Trim().Length -eq 0) {
            $next.Range.Delete()
        }
        break
    }
}

# --- 7. "Oppgaver til neste sprint" bullets ---
# Order matters: replace the (currently unique) "Lage Kortvelgesystem" bullet
# with its final text *before* introducing new "Lage Kortvelgesystem..." text
# via the "Designe ui" replacement below, to avoid Find matching the wrong run.
$d.Content.Find.Execute("Lage Kortvelgesystem", $true, $false, $false, $false, $false, $true, 1, $false, "Fullføre tekst", 2) | Out-Null
$d.Content.Find.Execute("Lage meny", $true, $false, $false, $false, $false, $true, 1, $false, "Fullføre klassediagram", 2) | Out-Null
$d.Content.Find.Execute("Designe ui", $true, $false, $false, $false, $false, $true, 1, $false, "Lage Kortvelgesystem som er GUI basert", 2) | Out-Null

# --- 8. Add a new bullet "Lage Presentasjon" after the last bullet ---
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$newp = $d.Paragraphs.Item($d.Paragraphs.Count)
$newp.Range.Text = "Lage Presentasjon"
